$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.408.86"
$ws.Range("E2").Value = "  +0.41%  "

$ws.Range("D3").Value = "2.647.39"
$ws.Range("E3").Value = "  +0.71%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.96%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "2.646.57"
$ws.Range("E9").Value = "  +0.71%  "

$ws.Range("E10").Value = "  +8.09%  "

$ws.Range("E11").Value = "  -0.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.34%  "

$ws.Range("E13").Value = "  +1.88%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.08%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000193"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.54%  "

$ws.Range("D16").Value = "3.127.00"
$ws.Range("E16").Value = "  +0.62%  "

$ws.Range("D17").Value = "68.309.89"
$ws.Range("E17").Value = "  +0.42%  "

$ws.Range("D18").Value = "2.640.39"
$ws.Range("E18").Value = "  -1.05%  "

$ws.Range("E19").Value = "  +1.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "365.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.67%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.78%  "

$ws.Range("E24").Value = "  +0.33%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.74%  "

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("E27").Value = "  -1.44%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000107"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.99%  "

$ws.Range("E29").Value = "  +0.65%  "

$ws.Range("E30").Value = "  -0.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "573.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.35%  "

$ws.Range("E34").Value = "  +0.93%  "

$ws.Range("E35").Value = "  +3.59%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.60"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.94"
$ws.Range("D38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.38%  "

$ws.Range("E40").Value = "  +0.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.375"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.51%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.89%  "

$ws.Range("D43").Value = "0.0₆0337"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.96%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.58%  "

$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "157.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.95%  "

$ws.Range("E49").Value = "  +1.96%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.45%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.95%  "
